$wb = $excel.ActiveWorkbook

# @@ -1159,22 +1159,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 186.55556
$ws.Range("I11").Value = 186.55556
$ws.Range("K11").Value = 186.55556
$ws.Range("M11").Value = -46.55556000000001

# @@ -1450,25 +1450,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 564.2
$ws.Range("J17").Value = 561.6667
$ws.Range("L17").Value = 1685.0001
$ws.Range("N17").Value = -2021.0001

# @@ -1551,25 +1551,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 223
$ws.Range("J19").Value = 346
$ws.Range("L19").Value = 346
$ws.Range("N19").Value = -696

# @@ -2194,22 +2194,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("L32").Value = 2000
$ws.Range("N32").Value = -2652

# @@ -2243,22 +2243,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 439.9
$ws.Range("I33").Value = 382.11765
$ws.Range("K33").Value = 382.11765
$ws.Range("M33").Value = -153.11765

# @@ -2647,25 +2647,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 693.3333
$ws.Range("I41").Value = 286.66666
$ws.Range("J41").Value = 1100
$ws.Range("K41").Value = 286.66666
$ws.Range("L41").Value = 1100
$ws.Range("M41").Value = 153.33334
$ws.Range("N41").Value = -1980

# @@ -3241,25 +3241,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3424.1
$ws.Range("I53").Value = 425
$ws.Range("J53").Value = 4173.875
$ws.Range("K53").Value = 425
$ws.Range("L53").Value = 4173.875
$ws.Range("M53").Value = 212
$ws.Range("N53").Value = -5447.875

# @@ -4095,22 +4095,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1379.2
$ws.Range("I70").Value = 1156
$ws.Range("K70").Value = 3468
$ws.Range("M70").Value = -3198

# @@ -4242,22 +4242,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1379.2
$ws.Range("I73").Value = 1156
$ws.Range("K73").Value = 3468
$ws.Range("M73").Value = -2532

# @@ -7187,25 +7187,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 152172.14
$ws.Range("I132").Value = 152172.14
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 456516.42
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -453986.42
$ws.Range("N132").ClearContents()

# @@ -7239,22 +7236,22 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 51746.668
$ws.Range("J133").Value = 51746.668
$ws.Range("L133").Value = 51746.668
$ws.Range("N133").Value = -61866.668

# @@ -7435,25 +7432,25 @@
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1594.3903
$ws.Range("I137").Value = 1472.6538
$ws.Range("J137").Value = 1805.4
$ws.Range("K137").Value = 4417.9614
$ws.Range("L137").Value = 5416.200000000001
$ws.Range("M137").Value = -1867.9614
$ws.Range("N137").Value = -10516.2

# @@ -7780,25 +7777,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3110.375
$ws.Range("J2").Value = 4753.25
$ws.Range("L2").Value = 4753.25
$ws.Range("N2").Value = -4979.25

# @@ -11290,22 +11287,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 796.37036
$ws.Range("I74").Value = 358.7143
$ws.Range("J74").Value = 1267.6923
$ws.Range("K74").Value = 358.7143
$ws.Range("L74").Value = 1267.6923
$ws.Range("M74").Value = 515.2857
$ws.Range("N74").Value = -3015.6923

# @@ -11434,22 +11434,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 796.37036
$ws.Range("I77").Value = 358.7143
$ws.Range("J77").Value = 1267.6923
$ws.Range("K77").Value = 1793.5715
$ws.Range("L77").Value = 6338.461499999999
$ws.Range("M77").Value = 2574.4285
$ws.Range("N77").Value = -15074.4615

# @@ -13327,25 +13330,25 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3110.375
$ws.Range("J116").Value = 4753.25
$ws.Range("L116").Value = 4753.25
$ws.Range("N116").Value = -9341.25

# @@ -14448,22 +14451,22 @@
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 49994
$ws.Range("J139").Value = 49994
$ws.Range("L139").Value = 49994
$ws.Range("N139").Value = -60274

# @@ -14735,25 +14738,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3110.375
$ws.Range("J3").Value = 4753.25
$ws.Range("L3").Value = 4753.25
$ws.Range("N3").Value = -4981.25

# @@ -18745,25 +18748,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1492.3043
$ws.Range("I86").Value = 1305.9375
$ws.Range("J86").Value = 1918.2858
$ws.Range("K86").Value = 1305.9375
$ws.Range("L86").Value = 1918.2858
$ws.Range("M86").Value = -182.9375
$ws.Range("N86").Value = -4164.2858

# @@ -18892,25 +18895,25 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1492.3043
$ws.Range("I89").Value = 1305.9375
$ws.Range("J89").Value = 1918.2858
$ws.Range("K89").Value = 6529.6875
$ws.Range("L89").Value = 9591.429
$ws.Range("M89").Value = -913.6875
$ws.Range("N89").Value = -20823.429

# @@ -21079,22 +21082,22 @@
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2815.6287
$ws.Range("I134").Value = 3002.9
$ws.Range("K134").Value = 9008.700000000001
$ws.Range("M134").Value = -6473.700000000001

# @@ -27385,22 +27388,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1236
$ws.Range("I122").Value = 1236
$ws.Range("K122").Value = 3708
$ws.Range("M122").Value = -1258

# @@ -27872,22 +27875,22 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3066.087
$ws.Range("I132").Value = 990.8333
$ws.Range("K132").Value = 2972.4999
$ws.Range("M132").Value = -442.4998999999998

# @@ -27973,25 +27976,25 @@
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 854.625
$ws.Range("I134").Value = 813.3333
$ws.Range("J134").Value = 978.5
$ws.Range("K134").Value = 2439.9999
$ws.Range("L134").Value = 2935.5
$ws.Range("M134").Value = 95.0001000000002
$ws.Range("N134").Value = -8005.5

# @@ -31780,25 +31783,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1213.1923
$ws.Range("J68").Value = 1221.72
$ws.Range("L68").Value = 3665.16
$ws.Range("N68").Value = -5287.16

# @@ -31936,25 +31939,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1213.1923
$ws.Range("J71").Value = 1221.72
$ws.Range("L71").Value = 10995.48
$ws.Range("N71").Value = -19107.48

# @@ -34057,25 +34060,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 466.2857
$ws.Range("I113").Value = 424
$ws.Range("J113").Value = 508.57144
$ws.Range("K113").Value = 1272
$ws.Range("L113").Value = 1525.71432
$ws.Range("M113").Value = 898
$ws.Range("N113").Value = -5865.71432

# @@ -34473,25 +34476,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 577
$ws.Range("J121").Value = 750
$ws.Range("L121").Value = 2250
$ws.Range("N121").Value = -4870

# @@ -34975,25 +34978,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 804.75
$ws.Range("J131").Value = 822.23956
$ws.Range("L131").Value = 2466.71868
$ws.Range("N131").Value = -12546.71868

# @@ -35027,25 +35030,25 @@
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1478.8
$ws.Range("I132").Value = 598.5
$ws.Range("J132").Value = 2065.6667
$ws.Range("K132").Value = 5386.5
$ws.Range("L132").Value = 18591.0003
$ws.Range("M132").Value = -2856.5
$ws.Range("N132").Value = -23651.0003

# @@ -41900,22 +41903,22 @@
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 27333.143
$ws.Range("I132").Value = 3714.1428
$ws.Range("K132").Value = 11142.4284
$ws.Range("M132").Value = -8612.428400000001

# @@ -42720,25 +42723,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3821.4546
$ws.Range("I7").Value = 4263.3335
$ws.Range("J7").Value = 2874.5715
$ws.Range("K7").Value = 4263.3335
$ws.Range("L7").Value = 2874.5715
$ws.Range("M7").Value = -4151.3335
$ws.Range("N7").Value = -3098.5715

# @@ -47725,22 +47728,22 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 35065.832
$ws.Range("J110").Value = 35065.832
$ws.Range("L110").Value = 35065.832
$ws.Range("N110").Value = -43245.832

# @@ -48494,25 +48497,25 @@
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3821.4546
$ws.Range("I126").Value = 4263.3335
$ws.Range("J126").Value = 2874.5715
$ws.Range("K126").Value = 12790.0005
$ws.Range("L126").Value = 8623.7145
$ws.Range("M126").Value = -10320.0005
$ws.Range("N126").Value = -13563.7145

# @@ -55724,25 +55727,25 @@
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1313.7333
$ws.Range("I132").Value = 1115.1613
$ws.Range("J132").Value = 1753.4286
$ws.Range("K132").Value = 3345.4839
$ws.Range("L132").Value = 5260.2858
$ws.Range("M132").Value = -815.4839000000002
$ws.Range("N132").Value = -10320.2858
